$d = $word.ActiveDocument

# Fix spelling: economicos -> econômicos
$d.Content.Find.Execute("economicos", $false, $false, $false, $false, $false,
                         $true, 1, $false, "econômicos", 2)

# Fix spelling: funcionario -> funcionário
$d.Content.Find.Execute("funcionario", $false, $false, $false, $false, $false,
                         $true, 1, $false, "funcionário", 2)
